$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.888.45"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3
$ws.Range("D3").Value = "3.576.79"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "201.25"
$ws.Range("E5").Value = "  +7.47%  "

# Row 6
$ws.Range("D6").Value = "566.35"
$ws.Range("E6").Value = "  -1.15%  "

# Row 7
$ws.Range("D7").Value = "3.571.29"
$ws.Range("E7").Value = "  +0.23%  "

# Row 8
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  -0.42%  "

# Row 9
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("D10").Value = "0.676"
$ws.Range("E10").Value = "  +0.87%  "

# Row 11
$ws.Range("D11").Value = "59.85"
$ws.Range("E11").Value = "  +8.05%  "

# Row 12
$ws.Range("E12").Value = "  -0.98%  "

# Row 13
$ws.Range("D13").Value = "0.0000279"
$ws.Range("E13").Value = "  +6.84%  "

# Row 14
$ws.Range("D14").Value = "10.19"
$ws.Range("E14").Value = "  +3.97%  "

# Row 15
$ws.Range("D15").Value = "4.150.80"
$ws.Range("E15").Value = "  +0.42%  "

# Row 16
$ws.Range("D16").Value = "3.575.80"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
$ws.Range("E17").Value = "  +0.89%  "

# Row 18
$ws.Range("D18").Value = "18.93"
$ws.Range("E18").Value = "  +3.49%  "

# Row 19
$ws.Range("D19").Value = "67.624.95"
$ws.Range("E19").Value = "  +1.47%  "

# Row 20
$ws.Range("D20").Value = "12.18"
$ws.Range("E20").Value = "  +0.87%  "

# Row 21
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("D22").Value = "401.16"
$ws.Range("E22").Value = "  +2.76%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "4.16"
$ws.Range("E23").Value = "  -1.41%  "

# Row 24
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  +11.73%  "

# Row 25
$ws.Range("D25").Value = "84.49"
$ws.Range("E25").Value = "  -1.12%  "

# Row 26
$ws.Range("D26").Value = "2.87"
$ws.Range("E26").Value = "  -1.80%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "3.88"
$ws.Range("E27").Value = "  +8.79%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "12.36"
$ws.Range("E28").Value = "  -0.39%  "

# Row 29
$ws.Range("E29").Value = "  +3.14%  "

# Row 30
$ws.Range("D30").Value = "7.68"
$ws.Range("E30").Value = "  +0.92%  "

# Row 31
$ws.Range("D31").Value = "31.39"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32
$ws.Range("D32").Value = "669.89"
$ws.Range("E32").Value = "  +6.64%  "

# Row 33
$ws.Range("D33").Value = "12.05"
$ws.Range("E33").Value = "  -1.14%  "

# Row 34
$ws.Range("D34").Value = "63.28"
$ws.Range("E34").Value = "  -0.21%  "

# Row 35
$ws.Range("D35").Value = "0.112"
$ws.Range("E35").Value = "  -1.70%  "

# Row 36
$ws.Range("D36").Value = "40.94"
$ws.Range("E36").Value = "  -2.22%  "

# Row 37
$ws.Range("D37").Value = "0.406"
$ws.Range("E37").Value = "  +0.70%  "

# Row 38
$ws.Range("E38").Value = "  -0.39%  "

# Row 39
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  +9.89%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41
$ws.Range("D41").Value = "3.187.50"
$ws.Range("E41").Value = "  +3.19%  "

# Row 42
$ws.Range("D42").Value = "0.132"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  +2.45%  "

# Row 45
$ws.Range("E45").Value = "  +12.50%  "

# Row 46
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.130"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.07"
$ws.Range("E48").Value = "  -1.27%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.61"
$ws.Range("E49").Value = "  +9.99%  "

# Row 50
$ws.Range("E50").Value = "  +2.18%  "

# Row 51
$ws.Range("D51").Value = "138.38"
$ws.Range("E51").Value = "  -0.57%  "
